$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Asistido Copilot" effort hours for the Bloque 2: Desarrollo
# CORE rows that previously had no value in column D.
$ws.Range("D19").Value = 2
$ws.Range("D20").Value = 2
$ws.Range("D21").Value = 3
$ws.Range("D22").Value = 3

# Move the view/selection to where the user was working (D23, the
# "Asistido Copilot" subtotal for that block) and scroll the window down
# so row 9 is at the top, matching the saved workbook view.
$win = $excel.ActiveWindow
$win.ScrollRow = 9
$win.ScrollColumn = 1
$ws.Range("D23").Select()
